$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range('D2').Value = '42.210.70'
$ws.Range('E2').Value = '  -0.39%  '
$ws.Range('D3').Value = '2.240.93'
$ws.Range('E3').Value = '  +0.09%  '
$ws.Range('E4').Value = '  +0.16%  '
$ws.Range('D5').Value = "'242.80"
$ws.Range('E5').Value = '  -1.21%  '
$ws.Range('D6').Value = "'0.626"
$ws.Range('E6').Value = '  -0.50%  '
$ws.Range('D7').Value = "'74.08"
$ws.Range('E7').Value = '  -0.49%  '
$ws.Range('E8').Value = '  +0.12%  '
$ws.Range('D9').Value = "'0.600"
$ws.Range('E9').Value = '  -3.44%  '
$ws.Range('D10').Value = "'42.28"
$ws.Range('E10').Value = '  -2.78%  '
$ws.Range('E11').Value = '  -0.02%  '
$ws.Range('B12').Value = 'TRON'
$ws.Range('C12').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('D12').Value = "'0.103"
$ws.Range('E12').Value = '  +0.30%  '
$ws.Range('B13').Value = 'Polkadot'
$ws.Range('C13').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D13').Value = "'6.93"
$ws.Range('E13').Value = '  -3.00%  '
$ws.Range('D14').Value = '2.574.24'
$ws.Range('E14').Value = '  +0.05%  '
$ws.Range('E15').Value = '  -1.06%  '
$ws.Range('E16').Value = '  -1.97%  '
$ws.Range('D17').Value = '2.239.82'
$ws.Range('E17').Value = '  +0.23%  '
$ws.Range('D18').Value = '42.106.89'
$ws.Range('E18').Value = '  -0.39%  '
$ws.Range('E19').Value = '  -5.32%  '
$ws.Range('E20').Value = '  +0.66%  '
$ws.Range('D21').Value = "'72.74"
$ws.Range('E21').Value = '  +0.89%  '
$ws.Range('D22').Value = "'11.31"
$ws.Range('E22').Value = '  +9.42%  '
$ws.Range('D23').Value = "'230.18"
$ws.Range('E23').Value = '  -0.66%  '
$ws.Range('E24').Value = '  -6.24%  '
$ws.Range('E25').Value = '  -0.04%  '
$ws.Range('D26').Value = "'11.42"
$ws.Range('E26').Value = '  -3.22%  '
$ws.Range('E27').Value = '  -0.39%  '
$ws.Range('D28').Value = "'2.28"
$ws.Range('E28').Value = '  -1.02%  '
$ws.Range('D30').Value = "'167.25"
$ws.Range('E30').Value = '  +0.18%  '
$ws.Range('E31').Value = '  -1.49%  '
$ws.Range('D32').Value = "'5.64"
$ws.Range('E32').Value = '  -4.49%  '
$ws.Range('D33').Value = "'0.0807"
$ws.Range('E33').Value = '  -0.89%  '
$ws.Range('D34').Value = "'29.61"
$ws.Range('E34').Value = '  -1.38%  '
$ws.Range('E35').Value = '  -0.49%  '
$ws.Range('D36').Value = "'0.111"
$ws.Range('E36').Value = '  -6.79%  '
$ws.Range('E37').Value = '  -5.05%  '
$ws.Range('E38').Value = '  -2.01%  '
$ws.Range('E39').Value = '  -1.75%  '
$ws.Range('E40').Value = '  -2.11%  '
$ws.Range('E41').Value = '  +1.20%  '
$ws.Range('D42').Value = "'64.64"
$ws.Range('E42').Value = '  +1.54%  '
$ws.Range('E43').Value = '  -1.49%  '
$ws.Range('E44').Value = '  -1.43%  '
$ws.Range('D45').Value = "'104.48"
$ws.Range('E45').Value = '  -1.35%  '
$ws.Range('D46').Value = "'0.101"
$ws.Range('E46').Value = '  -1.96%  '
$ws.Range('E47').Value = '  -0.50%  '
$ws.Range('E48').Value = '  -0.98%  '
$ws.Range('E49').Value = '  -2.49%  '
$ws.Range('E50').Value = '  -1.73%  '
$ws.Range('D51').Value = '2.448.28'
$ws.Range('E51').Value = '  -0.11%  '

# Reset style on cells where a leading apostrophe was used to force
# text interpretation, so no stray "quote prefix" cell style lingers.
$ws.Range('D5,D6,D7,D9,D10,D12,D13,D21,D22,D23,D26,D28,D30,D32,D33,D34,D36,D42,D45,D46').Style = "Normal"
